# Add cross-sectional Area columns alongside the existing Discharge (Q)
# columns, and a small side "total" summary block (J/K) that pulls the
# running Area/Q totals together - per commit message: "add area to Q
# files stn4".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers -------------------------------------------------
# G/H mirror the existing "Q"/"Qtotal" pattern (col E/F) but for Area.
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
# J/K is a small side summary: total area, total discharge.
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# Re-apply the existing segment-midpoint formula across D3:D8 so Excel
# collapses it back into a shared formula group (same text, same values).
$ws.Range("D3:D8").Formula = "=(A3/100+(A4/100-A3/100)/2)"

# --- Per-segment incremental area (trapezoid-ish strip: depth * width) --
# G2 uses 0 as the "previous" midpoint (mirrors D2 = A2/100 itself being
# the first midpoint), rows 3-15 use the previous row's midpoint.
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- Totals ---------------------------------------------------------------
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- View state: scroll right a bit and select the new summary cells ----
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("J2:K2").Select() | Out-Null
